$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.978.83"

$ws.Range("D3").Value = "1.870.02"
$ws.Range("E3").Value = "  -2.58%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.21"
$ws.Range("E5").Value = "  -3.31%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5048"
$ws.Range("E7").Value = "  -2.87%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3959"
$ws.Range("E8").Value = "  -3.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08207"
$ws.Range("E9").Value = "  -3.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.13"
$ws.Range("E10").Value = "  -2.57%  "

$ws.Range("E11").Value = "  -3.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.46"
$ws.Range("E12").Value = "  +4.84%  "

$ws.Range("D13").Value = "1.866.87"
$ws.Range("E13").Value = "  -2.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.293"
$ws.Range("E14").Value = "  -1.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.187"
$ws.Range("E15").Value = "  -2.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.86"
$ws.Range("E17").Value = "  -3.84%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001087"
$ws.Range("E18").Value = "  -2.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06424"
$ws.Range("E19").Value = "  -4.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.11"
$ws.Range("E20").Value = "  -0.75%  "

$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("D22").Value = "29.983.23"
$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.847"
$ws.Range("E23").Value = "  -2.90%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  -1.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.166"
$ws.Range("E25").Value = "  -2.52%  "

$ws.Range("D26").Value = "2.089.13"
$ws.Range("E26").Value = "  -2.42%  "

$ws.Range("E27").Value = "  +1.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.64"
$ws.Range("E28").Value = "  +0.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.211"
$ws.Range("E29").Value = "  -9.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.21"
$ws.Range("E30").Value = "  -1.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.072"
$ws.Range("E31").Value = "  -0.14%  "

$ws.Range("E32").Value = "  -1.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.939"
$ws.Range("E33").Value = "  -2.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.630"
$ws.Range("E34").Value = "  -0.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02438"
$ws.Range("E35").Value = "  -2.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.212"
$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06356"
$ws.Range("E37").Value = "  -3.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2140"
$ws.Range("E38").Value = "  -2.98%  "

$ws.Range("E39").Value = "  -4.84%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.480"
$ws.Range("E40").Value = "  -4.90%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.220"
$ws.Range("E41").Value = "  -2.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6309"
$ws.Range("E42").Value = "  -3.14%  "

$ws.Range("E43").Value = "  -3.11%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.96"
$ws.Range("E45").Value = "  -2.54%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5905"
$ws.Range("E46").Value = "  -4.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.089"
$ws.Range("E47").Value = "  +0.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.625"
$ws.Range("E48").Value = "  -3.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.58"
$ws.Range("E49").Value = "  -1.46%  "

$ws.Range("E50").Value = "  -3.48%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.46"
$ws.Range("E51").Value = "  -2.82%  "
